# Applies the "--save" flag to npm install commands on the "node" sheet,
# adds the missing install command for "moment", and widens column C
# to fit the longer command strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("node")

# express
$ws.Range("C2").Value = "npm install express --save"
# ioredis
$ws.Range("C3").Value = "npm install ioredis --save"
# mongoose
$ws.Range("C4").Value = "npm install mongoose --save"
# moment (previously had no install-command cell)
$ws.Range("C5").Value = "npm install moment --save"
# express-session
$ws.Range("C6").Value = "npm install express-session --save"
# connect-redis
$ws.Range("C7").Value = "npm install connect-redis --save"

# Widen column C so the longer install commands are fully visible.
$ws.Columns.Item(3).ColumnWidth = 29.7
